$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "17:05"
